$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 114
$ws1.Range("F3").Value = 7491
$ws1.Range("F5").Value = 6
$ws1.Range("F6").Value = 447
$ws1.Range("F7").Value = 4078
$ws1.Range("F9").Value = 573
$ws1.Range("F11").Value = 653
$ws1.Range("F12").Value = 136

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 7

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 114
$ws4.Range("F4").Value = 7491
$ws4.Range("F7").Value = 6
$ws4.Range("F8").Value = 447
$ws4.Range("F9").Value = 4078
$ws4.Range("F11").Value = 573
$ws4.Range("F13").Value = 653
$ws4.Range("F14").Value = 7
$ws4.Range("F15").Value = 136
